$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.246.89"
$ws.Range("D3").Value = "1.906.43"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.37"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5265"
$ws.Range("E7").Value = "  +1.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3815"
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07282"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.85"
$ws.Range("E10").Value = "  +3.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9026"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08166"
$ws.Range("E12").Value = "  -3.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.25"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("E14").Value = "  +1.28%  "
$ws.Range("D15").Value = "1.465.38"
$ws.Range("E15").Value = "  -23.13%  "
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008652"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("E18").Value = "  +1.42%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "27.284.93"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.118"
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.517"
$ws.Range("E23").Value = "  +1.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "150.01"
$ws.Range("E24").Value = "  +2.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.307"
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.740"
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "116.74"
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.848"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.852"
$ws.Range("E30").Value = "  -1.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09244"
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8306"
$ws.Range("E32").Value = "  +4.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05065"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.988"
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.347"
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.719"
$ws.Range("E37").Value = "  +5.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5817"
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.080"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.162"
$ws.Range("E41").Value = "  +1.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.600"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "117.03"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1522"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4934"
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.644"
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "39.00"
$ws.Range("E49").Value = "  +3.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06125"
$ws.Range("E50").Value = "  +2.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.52"
$ws.Range("E51").Value = "  +0.77%  "
